# Rename the "temperature_c" sheet to "temperature" and make it the
# active/selected sheet in the workbook (it was previously "genotype").

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("temperature_c")
$ws.Name = "temperature"

# Selecting/activating this sheet updates the workbook's activeTab and
# this sheet's tabSelected, while clearing tabSelected on the sheet that
# was previously active (genotype).
$ws.Activate()
